# Helper: set a cell as plain text (apostrophe-prefix forces text type,
# preventing Excel from auto-converting True/False/dates to other types).
function Set-Text($ws, $cellref, $text) {
    $ws.Range($cellref).Value = "'" + $text
}

function Add-Hlink($ws, $cellref, $url, $disp) {
    $ws.Hyperlinks.Add($ws.Range($cellref), $url, "", "", $disp)
}

$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item(1)
$wsZh = $wb.Worksheets.Item(2)
$wsDe = $wb.Worksheets.Item(3)

# ---- Update cell values: row 2 now holds the 4c308480 file, row 3 the 0ad0c369 file ----
# (previously it was the other way around); 0ad0c369 also gets a refreshed status.
# overview row2
Set-Text $wsOverview "A2" '4c308480-f3bc-4c02-bc75-8e1ff678093f.md'
Set-Text $wsOverview "B2" 'e2e\4c308480-f3bc-4c02-bc75-8e1ff678093f.md'
Set-Text $wsOverview "C2" '.md'
Set-Text $wsOverview "D2" ''
Set-Text $wsOverview "E2" 'Handed back: in sync with en-US'
Set-Text $wsOverview "F2" 'Handed back: in sync with en-US'
Set-Text $wsOverview "G2" '2016-09-01 15:07:28'

# overview row3
Set-Text $wsOverview "A3" '0ad0c369-1e2e-4808-9f44-53f2bda72a7c.md'
Set-Text $wsOverview "B3" 'e2e\0ad0c369-1e2e-4808-9f44-53f2bda72a7c.md'
Set-Text $wsOverview "C3" '.md'
Set-Text $wsOverview "D3" ''
Set-Text $wsOverview "E3" 'Ready for handoff'
Set-Text $wsOverview "F3" 'Ready for handoff'
Set-Text $wsOverview "G3" '2016-09-01 15:08:53'

# zh-cn row2
Set-Text $wsZh "A2" '4c308480-f3bc-4c02-bc75-8e1ff678093f.md'
Set-Text $wsZh "B2" '.md'
Set-Text $wsZh "C2" 'Handed back: in sync with en-US'
Set-Text $wsZh "D2" 'e2e'
Set-Text $wsZh "E2" 'ht'
Set-Text $wsZh "F2" 'False'
Set-Text $wsZh "G2" '4c308480-f3bc-4c02-bc75-8e1ff678093f.05c92128196b8ab187d24a42b7be9f6a43537871.zh-cn.xlf'
Set-Text $wsZh "H2" '2016-09-01 15:07:22'
Set-Text $wsZh "I2" '4c308480-f3bc-4c02-bc75-8e1ff678093f.md'
Set-Text $wsZh "J2" '4c308480-f3bc-4c02-bc75-8e1ff678093f.05c92128196b8ab187d24a42b7be9f6a43537871.zh-cn.xlf'
Set-Text $wsZh "K2" '2016-09-01 15:08:01'
Set-Text $wsZh "L2" ''
Set-Text $wsZh "M2" 'True'
Set-Text $wsZh "N2" ''
Set-Text $wsZh "O2" 'False'
Set-Text $wsZh "P2" ''

# zh-cn row3
Set-Text $wsZh "A3" '0ad0c369-1e2e-4808-9f44-53f2bda72a7c.md'
Set-Text $wsZh "B3" '.md'
Set-Text $wsZh "C3" 'Ready for handoff'
Set-Text $wsZh "D3" 'e2e'
Set-Text $wsZh "E3" 'ht'
Set-Text $wsZh "F3" 'False'
Set-Text $wsZh "G3" '0ad0c369-1e2e-4808-9f44-53f2bda72a7c.22014c441ecf50705aef5da5fafbdc282f856f9e.zh-cn.xlf'
Set-Text $wsZh "H3" '2016-09-01 15:08:48'
Set-Text $wsZh "I3" '0ad0c369-1e2e-4808-9f44-53f2bda72a7c.md'
Set-Text $wsZh "J3" '0ad0c369-1e2e-4808-9f44-53f2bda72a7c.22014c441ecf50705aef5da5fafbdc282f856f9e.zh-cn.xlf'
Set-Text $wsZh "K3" '2016-09-01 15:08:01'
Set-Text $wsZh "L3" ''
Set-Text $wsZh "M3" 'True'
Set-Text $wsZh "N3" ''
Set-Text $wsZh "O3" 'False'
Set-Text $wsZh "P3" 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3fc73fc81bf0923891a701fca259846a99493c17/e2e/0ad0c369-1e2e-4808-9f44-53f2bda72a7c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aa466a04f8063bc9756ff1b8de29c0dfede14298/e2e/0ad0c369-1e2e-4808-9f44-53f2bda72a7c.md.'

# de-de row2
Set-Text $wsDe "A2" '4c308480-f3bc-4c02-bc75-8e1ff678093f.md'
Set-Text $wsDe "B2" '.md'
Set-Text $wsDe "C2" 'Handed back: in sync with en-US'
Set-Text $wsDe "D2" 'e2e'
Set-Text $wsDe "E2" 'ht'
Set-Text $wsDe "F2" 'False'
Set-Text $wsDe "G2" '4c308480-f3bc-4c02-bc75-8e1ff678093f.05c92128196b8ab187d24a42b7be9f6a43537871.de-de.xlf'
Set-Text $wsDe "H2" '2016-09-01 15:07:28'
Set-Text $wsDe "I2" '4c308480-f3bc-4c02-bc75-8e1ff678093f.md'
Set-Text $wsDe "J2" '4c308480-f3bc-4c02-bc75-8e1ff678093f.05c92128196b8ab187d24a42b7be9f6a43537871.de-de.xlf'
Set-Text $wsDe "K2" '2016-09-01 15:08:27'
Set-Text $wsDe "L2" ''
Set-Text $wsDe "M2" 'True'
Set-Text $wsDe "N2" ''
Set-Text $wsDe "O2" 'False'
Set-Text $wsDe "P2" ''

# de-de row3
Set-Text $wsDe "A3" '0ad0c369-1e2e-4808-9f44-53f2bda72a7c.md'
Set-Text $wsDe "B3" '.md'
Set-Text $wsDe "C3" 'Ready for handoff'
Set-Text $wsDe "D3" 'e2e'
Set-Text $wsDe "E3" 'ht'
Set-Text $wsDe "F3" 'False'
Set-Text $wsDe "G3" '0ad0c369-1e2e-4808-9f44-53f2bda72a7c.22014c441ecf50705aef5da5fafbdc282f856f9e.de-de.xlf'
Set-Text $wsDe "H3" '2016-09-01 15:08:53'
Set-Text $wsDe "I3" '0ad0c369-1e2e-4808-9f44-53f2bda72a7c.md'
Set-Text $wsDe "J3" '0ad0c369-1e2e-4808-9f44-53f2bda72a7c.22014c441ecf50705aef5da5fafbdc282f856f9e.de-de.xlf'
Set-Text $wsDe "K3" '2016-09-01 15:08:27'
Set-Text $wsDe "L3" ''
Set-Text $wsDe "M3" 'True'
Set-Text $wsDe "N3" ''
Set-Text $wsDe "O3" 'False'
Set-Text $wsDe "P3" 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3fc73fc81bf0923891a701fca259846a99493c17/e2e/0ad0c369-1e2e-4808-9f44-53f2bda72a7c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aa466a04f8063bc9756ff1b8de29c0dfede14298/e2e/0ad0c369-1e2e-4808-9f44-53f2bda72a7c.md.'

# ---- Rebuild hyperlinks to match the swapped rows ----
$wsOverview.Hyperlinks.Delete()
$wsZh.Hyperlinks.Delete()
$wsDe.Hyperlinks.Delete()

# overview hyperlinks
Add-Hlink $wsOverview "B2" 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3fc73fc81bf0923891a701fca259846a99493c17/e2e/4c308480-f3bc-4c02-bc75-8e1ff678093f.md' 'e2e\4c308480-f3bc-4c02-bc75-8e1ff678093f.md'
Add-Hlink $wsOverview "B3" 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3fc73fc81bf0923891a701fca259846a99493c17/e2e/0ad0c369-1e2e-4808-9f44-53f2bda72a7c.md' 'e2e\0ad0c369-1e2e-4808-9f44-53f2bda72a7c.md'

# zh-cn hyperlinks
Add-Hlink $wsZh "A2" 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3fc73fc81bf0923891a701fca259846a99493c17/e2e/4c308480-f3bc-4c02-bc75-8e1ff678093f.md' '4c308480-f3bc-4c02-bc75-8e1ff678093f.md'
Add-Hlink $wsZh "I2" 'https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/449d026b4973b0c864dee4d364f7101fd4c9f7a1/e2e/4c308480-f3bc-4c02-bc75-8e1ff678093f.md' '4c308480-f3bc-4c02-bc75-8e1ff678093f.md'
Add-Hlink $wsZh "A3" 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3fc73fc81bf0923891a701fca259846a99493c17/e2e/0ad0c369-1e2e-4808-9f44-53f2bda72a7c.md' '0ad0c369-1e2e-4808-9f44-53f2bda72a7c.md'
Add-Hlink $wsZh "I3" 'https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/449d026b4973b0c864dee4d364f7101fd4c9f7a1/e2e/0ad0c369-1e2e-4808-9f44-53f2bda72a7c.md' '0ad0c369-1e2e-4808-9f44-53f2bda72a7c.md'

# de-de hyperlinks
Add-Hlink $wsDe "A2" 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3fc73fc81bf0923891a701fca259846a99493c17/e2e/4c308480-f3bc-4c02-bc75-8e1ff678093f.md' '4c308480-f3bc-4c02-bc75-8e1ff678093f.md'
Add-Hlink $wsDe "I2" 'https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/49d5cd636a2c2da301cc0d0d8f2b1984665c3b1e/e2e/4c308480-f3bc-4c02-bc75-8e1ff678093f.md' '4c308480-f3bc-4c02-bc75-8e1ff678093f.md'
Add-Hlink $wsDe "A3" 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3fc73fc81bf0923891a701fca259846a99493c17/e2e/0ad0c369-1e2e-4808-9f44-53f2bda72a7c.md' '0ad0c369-1e2e-4808-9f44-53f2bda72a7c.md'
Add-Hlink $wsDe "I3" 'https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/49d5cd636a2c2da301cc0d0d8f2b1984665c3b1e/e2e/0ad0c369-1e2e-4808-9f44-53f2bda72a7c.md' '0ad0c369-1e2e-4808-9f44-53f2bda72a7c.md'

# ---- Widen column P (Error Detail) on the zh-cn and de-de sheets ----
$wsZh.Columns.Item(16).ColumnWidth = 39.17
$wsDe.Columns.Item(16).ColumnWidth = 39.17
